# Apply the edit described by the diff:
#  - Column A values for rows 2..12 change from "CDF" to "CDF(e) " (trailing space preserved)
#  - The selected/active cell in the sheet view changes from G24 to E21
#
# (The sharedStrings reshuffling visible in the raw XML diff is just an
# internal side-effect of Excel's string table management when the "CDF"
# text is edited to "CDF(e) "; the other cells that pointed at "target" and
# "ppb" keep the same displayed values.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = "CDF(e) "
}

$ws.Range("E21").Select()
